$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6418.26
$ws.Range("I40").Value = 5037.394
$ws.Range("J40").Value = 9098.764999999999
$ws.Range("K40").Value = 5037.394
$ws.Range("L40").Value = 9098.764999999999
$ws.Range("M40").Value = -4862.394
$ws.Range("N40").Value = -9448.764999999999
$ws.Range("H107").Value = 1222.2727
$ws.Range("I107").Value = 1387.125
$ws.Range("K107").Value = 1387.125
$ws.Range("M107").Value = 532.875
$ws.Range("H137").Value = 4415.643
$ws.Range("I137").Value = 1694.2
$ws.Range("J137").Value = 5927.5557
$ws.Range("K137").Value = 5082.6
$ws.Range("L137").Value = 17782.6671
$ws.Range("M137").Value = -2532.6
$ws.Range("N137").Value = -22882.6671
$ws.Range("H138").Value = 3403.1428
$ws.Range("I138").Value = 955.5
$ws.Range("J138").Value = 6666.6665
$ws.Range("K138").Value = 2866.5
$ws.Range("L138").Value = 19999.9995
$ws.Range("M138").Value = 2273.5
$ws.Range("N138").Value = -30279.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 809.2941
$ws.Range("I2").Value = 654.9167
$ws.Range("J2").Value = 1179.8
$ws.Range("K2").Value = 654.9167
$ws.Range("L2").Value = 1179.8
$ws.Range("M2").Value = -541.9167
$ws.Range("N2").Value = -1405.8
$ws.Range("H116").Value = 809.2941
$ws.Range("I116").Value = 654.9167
$ws.Range("J116").Value = 1179.8
$ws.Range("K116").Value = 654.9167
$ws.Range("L116").Value = 1179.8
$ws.Range("M116").Value = 1639.0833
$ws.Range("N116").Value = -5767.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 809.2941
$ws.Range("I3").Value = 654.9167
$ws.Range("J3").Value = 1179.8
$ws.Range("K3").Value = 654.9167
$ws.Range("L3").Value = 1179.8
$ws.Range("M3").Value = -540.9167
$ws.Range("N3").Value = -1407.8
$ws.Range("H60").Value = 81199
$ws.Range("J60").Value = 81199
$ws.Range("L60").Value = 81199
$ws.Range("N60").Value = -82397
$ws.Range("H80").Value = 1466.6666
$ws.Range("J80").Value = 1861.6666
$ws.Range("L80").Value = 1861.6666
$ws.Range("N80").Value = -3857.6666
$ws.Range("H83").Value = 1466.6666
$ws.Range("J83").Value = 1861.6666
$ws.Range("L83").Value = 9308.333000000001
$ws.Range("N83").Value = -19292.333
$ws.Range("H134").Value = 4165
$ws.Range("I134").Value = 2247.5
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 6742.5
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = -4207.5
$ws.Range("N134").Value = -29070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2825
$ws.Range("I6").Value = 1766.6666
$ws.Range("J6").Value = 6000
$ws.Range("K6").Value = 1766.6666
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = -1653.6666
$ws.Range("N6").Value = -6226
$ws.Range("H7").Value = 6364.5293
$ws.Range("I7").Value = 8478.916999999999
$ws.Range("K7").Value = 8478.916999999999
$ws.Range("M7").Value = -8365.916999999999
$ws.Range("H22").Value = 716.6667
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -350
$ws.Range("N22").Value = -1450
$ws.Range("H31").Value = 7688.857
$ws.Range("I31").Value = 3533
$ws.Range("J31").Value = 9351.200000000001
$ws.Range("K31").Value = 3533
$ws.Range("L31").Value = 9351.200000000001
$ws.Range("M31").Value = -3238
$ws.Range("N31").Value = -9941.200000000001
$ws.Range("H34").Value = 7688.857
$ws.Range("I34").Value = 3533
$ws.Range("J34").Value = 9351.200000000001
$ws.Range("K34").Value = 3533
$ws.Range("L34").Value = 9351.200000000001
$ws.Range("M34").Value = -3331
$ws.Range("N34").Value = -9755.200000000001
$ws.Range("H132").Value = 2751.3845
$ws.Range("I132").Value = 2147.3333
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 6441.999899999999
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -3911.999899999999
$ws.Range("N132").Value = -35060
$ws.Range("H134").Value = 2402.1667
$ws.Range("I134").Value = 2082.6
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 6247.799999999999
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -3712.799999999999
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 58.142857
$ws.Range("I12").Value = 3.5
$ws.Range("J12").Value = 99.125
$ws.Range("K12").Value = 10.5
$ws.Range("L12").Value = 297.375
$ws.Range("M12").Value = 162.5
$ws.Range("N12").Value = -643.375
$ws.Range("H22").Value = 201
$ws.Range("J22").Value = 201
$ws.Range("L22").Value = 603
$ws.Range("N22").Value = -941
$ws.Range("H27").Value = 201
$ws.Range("J27").Value = 201
$ws.Range("L27").Value = 603
$ws.Range("N27").Value = -807
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents() | Out-Null
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents() | Out-Null
$ws.Range("H68").Value = 2455.3333
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2455.3333
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 7365.999899999999
$ws.Range("M68").ClearContents() | Out-Null
$ws.Range("N68").Value = -8987.999899999999
$ws.Range("H71").Value = 2455.3333
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2455.3333
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 22097.9997
$ws.Range("M71").ClearContents() | Out-Null
$ws.Range("N71").Value = -30209.9997
$ws.Range("H117").Value = 88.42856999999999
$ws.Range("J117").Value = 91.5
$ws.Range("L117").Value = 274.5
$ws.Range("N117").Value = -7158.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6764.5557
$ws.Range("I46").Value = 4333.6665
$ws.Range("J46").Value = 7980
$ws.Range("K46").Value = 4333.6665
$ws.Range("L46").Value = 7980
$ws.Range("M46").Value = -4145.6665
$ws.Range("N46").Value = -8356
$ws.Range("H132").Value = 3300
$ws.Range("I132").Value = 3300
$ws.Range("K132").Value = 9900
$ws.Range("M132").Value = -7370

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3923.6667
$ws.Range("I122").Value = 2372.8
$ws.Range("K122").Value = 7118.400000000001
$ws.Range("M122").Value = -4668.400000000001
$ws.Range("H132").Value = 3483.3333
$ws.Range("I132").Value = 3483.3333
$ws.Range("K132").Value = 10449.9999
$ws.Range("M132").Value = -7919.999899999999
